$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 157
$ws.Range("F7").Value = 5618
$ws.Range("C8").Value = '北京·AKB48 Team SH IDO SPECIAL 公演+握手会'
$ws.Range("D8").Value = '亦庄荣昌东街6号 北京亦创国际会展中心'
$ws.Range("E8").Value = '2024.08.03 12:30-08.03 17:30'
$ws.Range("F8").Value = 83
$ws.Range("G8").Value = 258
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=87794'
$ws.Range("I8").Value = '//i1.hdslb.com/bfs/openplatform/202406/nAWgPxcz1718941471491.png'
$ws.Range("C9").Value = '北京·GSCG动漫节'
$ws.Range("D9").Value = '石景山路68号 北京首钢会展中心'
$ws.Range("E9").Value = '2024.08.03 09:00-08.05 17:30'
$ws.Range("F9").Value = 7588
$ws.Range("G9").Value = 90
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=88048'
$ws.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202406/tBiF4jzt1719210068620.jpeg'
$ws.Range("C10").Value = '北京·一人之下爱好者聚会【免票活动】'
$ws.Range("D10").Value = '王府井大街88号 北京王府井银泰in88购物中心'
$ws.Range("E10").Value = '2024.08.03 14:00-08.03 18:00'
$ws.Range("F10").Value = 98
$ws.Range("G10").Value = 58
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=89933'
$ws.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202407/t8J245E61722321042395.jpeg'
$ws.Range("C11").Value = '北京·万游引力S8 知名配音演员 秦紫翼 内场见面签售会'
$ws.Range("E11").Value = '2024.08.03 11:00-08.03 17:00'
$ws.Range("F11").Value = 66
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=87543'
$ws.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202406/ovDCeTCx1718612781842.jpeg'
$ws.Range("C12").Value = '北京·万游引力S8 知名配音演员 续续点灯 内场见面签售会'
$ws.Range("E12").Value = '2024.08.03 12:00-08.03 17:00'
$ws.Range("F12").Value = 52
$ws.Range("G12").Value = 138
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=87542'
$ws.Range("I12").Value = '//i0.hdslb.com/bfs/openplatform/202406/c2yFURYC1718614591976.jpeg'
$ws.Range("C13").Value = '北京·万游引力国潮动漫嘉年华s8'
$ws.Range("D13").Value = '金蝉西路甲1号 北京酷车国际汇展中心'
$ws.Range("E13").Value = '2024.08.03 10:00-08.04 17:00'
$ws.Range("F13").Value = 3832
$ws.Range("G13").Value = 75
$ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=83930'
$ws.Range("I13").Value = '//i1.hdslb.com/bfs/openplatform/202406/0kySwWBG1718096478563.jpeg'
$ws.Range("F16").Value = 196
$ws.Range("F19").Value = 103
$ws.Range("F21").Value = 596
$ws.Range("F22").Value = 3875
$ws.Range("F23").Value = 131
$ws.Range("F25").Value = 5273
$ws.Range("F26").Value = 439
$ws.Range("F27").Value = 2085
$ws.Range("F28").Value = 131
$ws.Range("F29").Value = 347
$ws.Range("F30").Value = 7824
$ws.Range("F34").Value = 2171
$ws.Range("F36").Value = 1191
$ws.Range("F39").Value = 265
$ws.Range("F42").Value = 1176
$ws.Range("F43").Value = 1174
$ws.Range("F44").Value = 31
$ws.Range("F45").Value = 1325
$ws.Range("F46").Value = 2064
$ws.Range("F47").Value = 124
$ws.Range("F48").Value = 219
$ws.Range("F49").Value = 1217
$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 122
$ws.Range("F20").Value = 5
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 565
$ws.Range("F3").Value = 742
$ws.Range("F4").Value = 64
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 157
$ws.Range("F5").Value = 565
$ws.Range("F6").Value = 742
$ws.Range("F8").Value = 5618
$ws.Range("F9").Value = 7588
$ws.Range("C10").Value = '北京·一人之下爱好者聚会【免票活动】'
$ws.Range("F10").Value = 98
$ws.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202407/t8J245E61722321042395.jpeg'
$ws.Range("F11").Value = 3832
$ws.Range("F14").Value = 196
$ws.Range("F17").Value = 103
$ws.Range("F20").Value = 596
$ws.Range("F21").Value = 3875
$ws.Range("F23").Value = 131
$ws.Range("F25").Value = 5273
$ws.Range("F26").Value = 439
$ws.Range("F27").Value = 2085
$ws.Range("F28").Value = 131
$ws.Range("F29").Value = 347
$ws.Range("F30").Value = 7824
$ws.Range("F34").Value = 2172
$ws.Range("F36").Value = 1191
$ws.Range("F37").Value = 265
$ws.Range("F40").Value = 1176
$ws.Range("F41").Value = 1174
$ws.Range("F42").Value = 31
$ws.Range("F43").Value = 1325
$ws.Range("F44").Value = 2064
$ws.Range("F45").Value = 124
$ws.Range("F47").Value = 219
$ws.Range("F48").Value = 5
$ws.Range("F49").Value = 1217
